# Applies the cryptos-list refresh described in the commit
# "Updated cryptos list on Thu Jan 25 08:31:08 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.131.56"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.224.89"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'292.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'87.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "'30.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'0.0781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("D13").Value = "'6.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "2.571.91"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'13.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "2.233.28"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "'0.729"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "40.093.65"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0886"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'11.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.85%  "
$ws.Range("D21").Value = "'5.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "'65.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'237.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "'1.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "'22.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "'9.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'156.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "'31.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'4.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "'2.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.85%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'15.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("D39").Value = "'0.0985"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "'1.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "2.114.72"
$ws.Range("E41").Value = "  +8.15%  "
$ws.Range("D42").Value = "'3.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.21%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'17.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.10%  "
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "2.438.79"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'69.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.53%  "
